$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Header labels (team/stack names) ---
$ws.Range("A1").Value = "Toronto Blue Jays lefties"
$ws.Range("F1").Value = "Tampa Bay Rays righties"
$ws.Range("K1").Value = "Minnesota Twins lefties"

# --- Block A: Toronto Blue Jays lefties (columns A-D) ---
$ws.Range("A3").Value = "Solarte"
$ws.Range("B3").Value = 3800
$ws.Range("C3").Value = 6

$ws.Range("A4").Value = "Smoak"
$ws.Range("B4").Value = 3500
$ws.Range("C4").Value = 3

$ws.Range("A5").Value = "Morales"
$ws.Range("B5").Value = 2600
$ws.Range("C5").Value = 0

# Row 6 in block A has no player this time around
$ws.Range("A6:C6").ClearContents() | Out-Null

# --- Block F: Tampa Bay Rays righties (columns F-I) ---
$ws.Range("F3").Value = "Span"
$ws.Range("G3").Value = 3000
$ws.Range("H3").Value = 0

$ws.Range("F4").Value = "Cron"
$ws.Range("G4").Value = 2900
$ws.Range("H4").Value = 6

$ws.Range("F5").Value = "Duffy"
$ws.Range("G5").Value = 2800
$ws.Range("H5").Value = 6

$ws.Range("F6").Value = "Ramos"
$ws.Range("G6").Value = 2700
$ws.Range("H6").Value = 6

# --- Block K: Minnesota Twins lefties (columns K-N) ---
$ws.Range("K3").Value = "Mauer"
$ws.Range("L3").Value = 3200
$ws.Range("M3").Value = 6.2

$ws.Range("K4").Value = "Kepler"
$ws.Range("L4").Value = 3300
$ws.Range("M4").Value = 3.5

$ws.Range("K5").Value = "Rosario"
$ws.Range("L5").Value = 3800
$ws.Range("M5").Value = 34.4

$ws.Range("K6").Value = "Morrison "
$ws.Range("L6").Value = 3000
$ws.Range("M6").Value = 13

# --- Result labels under each block ---
$ws.Range("D8").Value = "Failure"
$ws.Range("I8").Value = "Failure"
$ws.Range("N8").Value = "Success"

# --- Selection shown in the saved sheet view ---
$ws.Range("N9").Select() | Out-Null
